# Adds EGID / EWID / STEUERBARESEINKOMMEN / AMOUNT columns to the FAKE_DATA
# sheet, inserting them in amongst the existing VERMÖGEN / EL-BEZUG / SH-BEZUG
# columns (N:P -> N:T) and filling in the sample values for rows 2-4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlShiftToRight = -4161

# Step 1: insert 2 new columns at N:O (pushes old VERMÖGEN/EL-BEZUG/SH-BEZUG from N,O,P to P,Q,R)
$ws.Range("N1:O4").Insert($xlShiftToRight)

# Step 2: insert 1 new column at Q (between the shifted VERMÖGEN (P) and EL-BEZUG (Q)), pushing EL-BEZUG/SH-BEZUG to R,S
$ws.Range("Q1:Q4").Insert($xlShiftToRight)

# Step 3: headers
$ws.Range("N1").Value = "EGID"
$ws.Range("O1").Value = "EWID"
$ws.Range("Q1").Value = "STEUERBARESEINKOMMEN"

# T1 is a brand new trailing header cell; copy the bold header formatting from
# its neighbour (S1) before writing the label, matching the other header cells
$ws.Range("S1").Copy($ws.Range("T1"))
$ws.Range("T1").Value = "AMOUNT"

# Step 4: data values
$ws.Range("N2").Value = 11111
$ws.Range("Q2").Value = 12345
$ws.Range("T2").Value = 100

$ws.Range("N3").Value = 22
$ws.Range("Q3").Value = 99999
$ws.Range("T3").Value = 2000

$ws.Range("N4").Value = 9
$ws.Range("O4").Value = 122
$ws.Range("Q4").Value = 400000
$ws.Range("T4").Value = -50

# Selection / active cell, matches the committed workbook view
$ws.Range("T5").Select()
